$d = $word.ActiveDocument

# Locate the word "Install" that precedes "required packages".
$r = $d.Content
$found = $r.Find.Execute("Install required packages", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)

$start = $r.Start

# Replace "Install" with "Load" (leaves " required packages" immediately after).
$rWord = $d.Range($start, $start + 7)
$rWord.Text = "Load"

# Force the trailing text to split into two distinct runs: "Load" and
# " required packages" -- toggling a character property and restoring it
# causes the engine to re-emit the touched span as its own run without
# altering the visible formatting.
$loadEnd = $start + 4
$restEnd = $loadEnd + 18
$afterEnd = $restEnd + 35

$rLoad = $d.Range($start, $loadEnd)
$rLoad.Bold = 1
$rLoad.Bold = 0

$rRest = $d.Range($loadEnd, $restEnd)
$rRest.Bold = 1
$rRest.Bold = 0

# Re-establish the pre-existing run boundary right after "...and read “" so
# that the following (untouched) "Climate station data.xlsx”" text does not
# stay fused to our edit.
$rQuote = $d.Range($restEnd, $afterEnd)
$rQuote.Bold = 1
$rQuote.Bold = 0
